# Fruta / hortaliza, semanal
# Insert a new weekly price row for "Feria Lagunitas de Puerto Montt - Piña"
# at row 90 (pushing the existing rows 90..143 down to 91..144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 90..143 down by one to make room for the new record.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly record.
$ws.Range("A90").Value = 4
$ws.Range("B90").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C90").Value = "Los Lagos"
$ws.Range("D90").Value = 44488
$ws.Range("E90").Value = 10
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100108
$ws.Range("H90").Value = "Tropicales y subtropicales"
$ws.Range("I90").Value = 100108005
$ws.Range("J90").Value = "Piña"
$ws.Range("K90").Value = "Caramelo"
$ws.Range("L90").Value = "Segunda"
$ws.Range("M90").Value = 200
$ws.Range("N90").Value = 25000
$ws.Range("O90").Value = 25500
$ws.Range("P90").Value = 25250
$ws.Range("Q90").Value = "$/caja 14 unidades"
$ws.Range("R90").Value = "Ecuador"
$ws.Range("S90").Value = 1804
$ws.Range("T90").Value = 14
